# Generate Report for Handoff
# Updates the localization-status workbook with a new handoff round:
#  - UUID token changes from 9e60c5a1-8274-476e-8f06-e24b070f381c
#    to d93a9d40-8581-46d8-8641-19ac82c3f509
#  - xliff content hash changes from af5a253d1d1f4e4a1e0e855c186784c35466088e
#    to c6ce4e430e783b9f8d29f3cb297000026aa2ad86
#  - a few timestamps advance

$wb = $excel.ActiveWorkbook

$oldUuid = "9e60c5a1-8274-476e-8f06-e24b070f381c"
$newUuid = "d93a9d40-8581-46d8-8641-19ac82c3f509"
$oldHash = "af5a253d1d1f4e4a1e0e855c186784c35466088e"
$newHash = "c6ce4e430e783b9f8d29f3cb297000026aa2ad86"

# All three "source file" hyperlinks point at the same GitHub blob URL
# (the commit sha portion of the URL is untouched by this edit).
$commitSha = "3fdcb2b19cb4d44c5e597bb8832e5f0ac5e1d9a9"
$sourceAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/" + $commitSha + "/e2e/" + $oldUuid + ".md"

# ---- Sheet "Overview" ----
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = ($newUuid + ".md")
$wsOverview.Range("G2").Value = "2016-08-27 16:56:10"

# Update the hyperlink display text on B2 while keeping the same target
# address/relationship. Re-creating the hyperlink (delete + add with the
# original address) is required because editing TextToDisplay in place
# duplicates the hyperlink entry in this runtime.
$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $sourceAddress, "", "", ("e2e\" + $newUuid + ".md")) | Out-Null

# ---- Sheet "zh-cn" ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $sourceAddress, "", "", ($newUuid + ".md")) | Out-Null
$wsZhCn.Range("A2").Value = ($newUuid + ".md")

$wsZhCn.Range("G2").Value = ($newUuid + "." + $newHash + ".zh-cn.xlf")
$wsZhCn.Range("H2").Value = "2016-08-27 16:56:05"

# ---- Sheet "de-de" ----
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $sourceAddress, "", "", ($newUuid + ".md")) | Out-Null
$wsDeDe.Range("A2").Value = ($newUuid + ".md")

$wsDeDe.Range("G2").Value = ($newUuid + "." + $newHash + ".de-de.xlf")
$wsDeDe.Range("H2").Value = "2016-08-27 16:56:10"
